# Update publish metadata (version bump + republish date) on the "Metadata" sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.0.0"
$meta.Range("B8").Value = "2025-06-05T14:31:57+02:00"

# Fix off app-4 error on no-basis-appointment: clear the stray "Note that FHIR
# strings SHALL NOT exceed 1MB in size" Comments text that had been duplicated
# onto the code/display children of Extension.value[x], and clear the
# erroneous "ele-1" Condition(s) entries that had leaked onto rows that
# shouldn't carry that invariant. Also normalise the RIM mapping casing for
# Extension.value[x].extension to lowercase "n/a".
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("N11").ClearContents()
$elements.Range("N12").ClearContents()

$elements.Range("AI4").ClearContents()
$elements.Range("AI6").ClearContents()
$elements.Range("AI8").ClearContents()
$elements.Range("AI9").ClearContents()
$elements.Range("AI10").ClearContents()
$elements.Range("AI11").ClearContents()
$elements.Range("AI12").ClearContents()
$elements.Range("AI13").ClearContents()

$elements.Range("AK8").Value = "n/a"
